# Append two new weekly blocks (week 31 / Hugo, week 32 / Laura) below the
# existing data (header + week 30 / Hugo + Laura) on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: week number
$ws.Range("A5").Value = 31

# Row 6: Peltier / Hugo attendance for week 31
$ws.Range("A6").Value = "Peltier"
$ws.Range("B6").Value = "Hugo"
$ws.Range("C6").Value = "No"
$ws.Range("D6").Value = "Yes"
$ws.Range("E6").Value = "Yes"
$ws.Range("F6").Value = "No"
$ws.Range("G6").Value = "No"

# Row 7: week number
$ws.Range("A7").Value = 32

# Row 8: Peltier / Laura attendance for week 32
$ws.Range("A8").Value = "Peltier"
$ws.Range("B8").Value = "Laura"
$ws.Range("C8").Value = "No"
$ws.Range("D8").Value = "Yes"
$ws.Range("E8").Value = "Holiday"
$ws.Range("F8").Value = "Yes"
$ws.Range("G8").Value = "No"
